# finalise uk update of data and parameters
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UK")

# Updated regression coefficients for B2:B5 (also picks up a new direct
# number-format style, same as the source workbook's diff: B2:B5 move off
# the default style index 0 onto a new style index 1).
$ws.Range("B2").Value = 0.5511963385465366
$ws.Range("B3").Value = 0.56378226458419012
$ws.Range("B4").Value = 0.3588496779494994
$ws.Range("B5").Value = 0.37499140584966989
$ws.Range("B2:B5").NumberFormat = "General"

# Move the active selection to D5 (was A19).
$ws.Range("D5").Select()
